# Resolving version conflicts in Excel documents (I hate merge conflicts).
#
# A new "description" tag-match block (h4 / class / prod-title) is inserted
# right after the existing "description -> a" block, which in turn shifts
# every subsequent block (the three "price" blocks) down by one 4-row group
# (i.e. by 5 rows, including the blank separator row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New block: description / h4 / class / prod-title (rows 16-18) -------
# Row 15 (Item_Attribute_Name / description) is unchanged.
$ws.Range("B16").Value = "h4"
$ws.Range("A17").Value = "html_tag_attribute_name"
$ws.Range("B17").Value = "class"
$ws.Range("A18").Value = "html_tag_attribute_value"
$ws.Range("B18").Value = "prod-title"

# --- Old "description / a" block moves down to rows 20-23 -----------------
$ws.Range("A20").Value = "Item_Attribute_Name"
$ws.Range("B20").Value = "description"
$ws.Range("A21").Value = "html_tag"
$ws.Range("B21").Value = "a"
$ws.Range("A22").Value = "html_tag_attribute_name"
$ws.Range("B22").ClearContents()
$ws.Range("A23").Value = "html_tag_attribute_value"
$ws.Range("B23").ClearContents()

# --- First "price" block moves down to rows 25-28 --------------------------
$ws.Range("A25").Value = "Item_Attribute_Name"
$ws.Range("B25").Value = "price"
$ws.Range("A26").Value = "html_tag"
$ws.Range("B26").Value = "div"
$ws.Range("A27").Value = "html_tag_attribute_name"
$ws.Range("B27").Value = "class"
$ws.Range("A28").Value = "html_tag_attribute_value"
$ws.Range("B28").Value = "price-wrapper price-medium"

# --- Second "price" block moves down to rows 30-33 --------------------------
$ws.Range("A30").Value = "Item_Attribute_Name"
$ws.Range("B30").Value = "price"
$ws.Range("A31").Value = "html_tag"
$ws.Range("B31").Value = "div"
$ws.Range("A32").Value = "html_tag_attribute_name"
$ws.Range("B32").Value = "class"
$ws.Range("A33").Value = "html_tag_attribute_value"
$ws.Range("B33").Value = "prodprice"

# --- Third "price" block is brand new, rows 35-38 ---------------------------
$ws.Range("A35").Value = "Item_Attribute_Name"
$ws.Range("B35").Value = "price"
$ws.Range("A36").Value = "html_tag"
$ws.Range("B36").Value = "span"
$ws.Range("A37").Value = "html_tag_attribute_name"
$ws.Range("B37").Value = "class"
$ws.Range("A38").Value = "html_tag_attribute_value"
$ws.Range("B38").Value = "amount"

# Update the active selection to match the authored edit (B18).
$ws.Range("B18").Select()
